$d = $word.ActiveDocument

$d.Content.Find.Execute("50×12=", $true, $false, $false, $false, $false, $true, 1, $false, "91×67=", 2)
$d.Content.Find.Execute("19×39=", $true, $false, $false, $false, $false, $true, 1, $false, "56×20=", 2)
$d.Content.Find.Execute("51×58=", $true, $false, $false, $false, $false, $true, 1, $false, "80×73=", 2)
$d.Content.Find.Execute("35×37=", $true, $false, $false, $false, $false, $true, 1, $false, "97×18=", 2)
$d.Content.Find.Execute("17×66=", $true, $false, $false, $false, $false, $true, 1, $false, "34×51=", 2)
$d.Content.Find.Execute("46×53=", $true, $false, $false, $false, $false, $true, 1, $false, "13×40=", 2)
$d.Content.Find.Execute("87×60=", $true, $false, $false, $false, $false, $true, 1, $false, "76×89=", 2)
$d.Content.Find.Execute("12×54=", $true, $false, $false, $false, $false, $true, 1, $false, "56×86=", 2)
$d.Content.Find.Execute("39×32=", $true, $false, $false, $false, $false, $true, 1, $false, "27×27=", 2)
$d.Content.Find.Execute("80×37=", $true, $false, $false, $false, $false, $true, 1, $false, "16×29=", 2)
$d.Content.Find.Execute("72×32=", $true, $false, $false, $false, $false, $true, 1, $false, "59×48=", 2)
$d.Content.Find.Execute("91×20=", $true, $false, $false, $false, $false, $true, 1, $false, "51×14=", 2)
$d.Content.Find.Execute("64×85=", $true, $false, $false, $false, $false, $true, 1, $false, "74×25=", 2)
$d.Content.Find.Execute("43×36=", $true, $false, $false, $false, $false, $true, 1, $false, "56×75=", 2)
$d.Content.Find.Execute("59×89=", $true, $false, $false, $false, $false, $true, 1, $false, "95×16=", 2)
$d.Content.Find.Execute("90×35=", $true, $false, $false, $false, $false, $true, 1, $false, "31×53=", 2)
$d.Content.Find.Execute("48×16=", $true, $false, $false, $false, $false, $true, 1, $false, "14×69=", 2)
$d.Content.Find.Execute("66×30=", $true, $false, $false, $false, $false, $true, 1, $false, "99×97=", 2)
$d.Content.Find.Execute("74×34=", $true, $false, $false, $false, $false, $true, 1, $false, "51×62=", 2)
$d.Content.Find.Execute("31×31=", $true, $false, $false, $false, $false, $true, 1, $false, "50×85=", 2)
$d.Content.Find.Execute("84×77=", $true, $false, $false, $false, $false, $true, 1, $false, "76×14=", 2)
$d.Content.Find.Execute("22×99=", $true, $false, $false, $false, $false, $true, 1, $false, "54×97=", 2)
$d.Content.Find.Execute("92×21=", $true, $false, $false, $false, $false, $true, 1, $false, "37×38=", 2)
$d.Content.Find.Execute("78×93=", $true, $false, $false, $false, $false, $true, 1, $false, "60×23=", 2)
$d.Content.Find.Execute("27×68=", $true, $false, $false, $false, $false, $true, 1, $false, "12×20=", 2)
